# Auto-generated edit script: refresh market-price-derived columns (H-N)
# on the Leve profit tables, per the scheduled-runner update.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 9596.263000000001
$ws.Range("I40").Value = 8314.182000000001
$ws.Range("J40").Value = 10118.593
$ws.Range("K40").Value = 8314.182000000001
$ws.Range("L40").Value = 10118.593
$ws.Range("M40").Value = -8139.182000000001
$ws.Range("N40").Value = -10468.593
$ws.Range("H106").Value = 15932.125
$ws.Range("I106").Value = 8701.75
$ws.Range("J106").Value = 37623.25
$ws.Range("K106").Value = 8701.75
$ws.Range("L106").Value = 37623.25
$ws.Range("M106").Value = -8070.75
$ws.Range("N106").Value = -38885.25
$ws.Range("H129").Value = 2461
$ws.Range("I129").Value = 1811.625
$ws.Range("K129").Value = 5434.875
$ws.Range("M129").Value = -434.875
$ws.Range("H138").Value = 4380.7646
$ws.Range("I138").Value = 3946.3845
$ws.Range("J138").Value = 4649.6665
$ws.Range("K138").Value = 11839.1535
$ws.Range("L138").Value = 13948.9995
$ws.Range("M138").Value = -6699.1535
$ws.Range("N138").Value = -24228.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2928
$ws.Range("I45").Value = 1923.2
$ws.Range("K45").Value = 1923.2
$ws.Range("M45").Value = -1546.2
$ws.Range("H63").Value = 5244.9375
$ws.Range("J63").Value = 8877.143
$ws.Range("L63").Value = 8877.143
$ws.Range("N63").Value = -10249.143
$ws.Range("H66").Value = 5244.9375
$ws.Range("J66").Value = 8877.143
$ws.Range("L66").Value = 44385.715
$ws.Range("N66").Value = -51249.715
$ws.Range("H97").Value = 2683.8928
$ws.Range("I97").Value = 2149.9565
$ws.Range("J97").Value = 5140
$ws.Range("K97").Value = 2149.9565
$ws.Range("L97").Value = 5140
$ws.Range("M97").Value = -1653.9565
$ws.Range("N97").Value = -6132
$ws.Range("H132").Value = 9374.315000000001
$ws.Range("I132").Value = 6197.375
$ws.Range("J132").Value = 11684.818
$ws.Range("K132").Value = 18592.125
$ws.Range("L132").Value = 35054.454
$ws.Range("M132").Value = -16062.125
$ws.Range("N132").Value = -40114.454

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 26944.299
$ws.Range("I31").Value = 2807.5
$ws.Range("K31").Value = 2807.5
$ws.Range("M31").Value = -2512.5
$ws.Range("H34").Value = 26944.299
$ws.Range("I34").Value = 2807.5
$ws.Range("K34").Value = 2807.5
$ws.Range("M34").Value = -2605.5
$ws.Range("H88").Value = 12611.333
$ws.Range("J88").Value = 12611.333
$ws.Range("L88").Value = 12611.333
$ws.Range("N88").Value = -13423.333
$ws.Range("H91").Value = 12611.333
$ws.Range("J91").Value = 12611.333
$ws.Range("L91").Value = 12611.333
$ws.Range("N91").Value = -15419.333
$ws.Range("H107").Value = 1412.35
$ws.Range("J107").Value = 2908
$ws.Range("L107").Value = 2908
$ws.Range("N107").Value = -6748
$ws.Range("H132").Value = 4147.4
$ws.Range("I132").Value = 2899.6875
$ws.Range("K132").Value = 8699.0625
$ws.Range("M132").Value = -6169.0625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 23211.125
$ws.Range("I2").Value = 1935.875
$ws.Range("J2").Value = 44486.375
$ws.Range("K2").Value = 11615.25
$ws.Range("L2").Value = 266918.25
$ws.Range("M2").Value = -11502.25
$ws.Range("N2").Value = -267144.25
$ws.Range("H62").Value = 3149.5
$ws.Range("J62").Value = 4002
$ws.Range("L62").Value = 12006
$ws.Range("N62").Value = -13378
$ws.Range("H65").Value = 3149.5
$ws.Range("J65").Value = 4002
$ws.Range("L65").Value = 36018
$ws.Range("N65").Value = -42882
$ws.Range("H113").Value = 1379
$ws.Range("I113").Value = 1070.4286
$ws.Range("J113").Value = 1533.2858
$ws.Range("K113").Value = 3211.2858
$ws.Range("L113").Value = 4599.857400000001
$ws.Range("M113").Value = -1041.2858
$ws.Range("N113").Value = -8939.857400000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 18503
$ws.Range("J21").Value = 20000
$ws.Range("L21").Value = 20000
$ws.Range("N21").Value = -20346
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("H30").Value = 18503
$ws.Range("J30").Value = 20000
$ws.Range("L30").Value = 20000
$ws.Range("N30").Value = -20210
$ws.Range("H80").Value = 3748.423
$ws.Range("I80").Value = 2847.389
$ws.Range("J80").Value = 5775.75
$ws.Range("K80").Value = 2847.389
$ws.Range("L80").Value = 5775.75
$ws.Range("M80").Value = -1849.389
$ws.Range("N80").Value = -7771.75
$ws.Range("H83").Value = 3748.423
$ws.Range("I83").Value = 2847.389
$ws.Range("J83").Value = 5775.75
$ws.Range("K83").Value = 14236.945
$ws.Range("L83").Value = 28878.75
$ws.Range("M83").Value = -9244.945
$ws.Range("N83").Value = -38862.75
$ws.Range("H113").Value = 7279.5454
$ws.Range("I113").Value = 1725
$ws.Range("J113").Value = 9362.5
$ws.Range("K113").Value = 1725
$ws.Range("L113").Value = 9362.5
$ws.Range("M113").Value = 445
$ws.Range("N113").Value = -13702.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1598.9286
$ws.Range("I16").Value = 1507.9166
$ws.Range("K16").Value = 1507.9166
$ws.Range("M16").Value = -1337.9166
$ws.Range("H122").Value = 8767.111000000001
$ws.Range("J122").Value = 11800.6
$ws.Range("L122").Value = 35401.8
$ws.Range("N122").Value = -40301.8
$ws.Range("H132").Value = 6358.5557
$ws.Range("I132").Value = 4194.9
$ws.Range("J132").Value = 9063.125
$ws.Range("K132").Value = 12584.7
$ws.Range("L132").Value = 27189.375
$ws.Range("M132").Value = -10054.7
$ws.Range("N132").Value = -32249.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 8945.385
$ws.Range("J13").Value = 8921.111000000001
$ws.Range("L13").Value = 8921.111000000001
$ws.Range("N13").Value = -9201.111000000001
$ws.Range("H122").Value = 2338.7273
$ws.Range("I122").Value = 1535.6842
$ws.Range("K122").Value = 4607.0526
$ws.Range("M122").Value = -2157.0526
$ws.Range("H132").Value = 9901.857
$ws.Range("I132").Value = 4745.9375
$ws.Range("K132").Value = 14237.8125
$ws.Range("M132").Value = -11707.8125
$ws.Range("H136").Value = 3250.1333
$ws.Range("I136").Value = 2872.2307
$ws.Range("J136").Value = 5706.5
$ws.Range("K136").Value = 8616.6921
$ws.Range("L136").Value = 17119.5
$ws.Range("M136").Value = -6066.6921
$ws.Range("N136").Value = -22219.5

